$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"   = -8.183
    "D14"  = -7.804
    "D16"  = -8.353
    "D21"  = -8.300999999999998
    "D23"  = -7.980000000000001
    "D25"  = -7.842999999999999
    "D26"  = -7.399000000000001
    "D29"  = -7.227000000000001
    "D40"  = -8.382
    "D53"  = -7.947999999999999
    "D57"  = -8
    "D59"  = -8.137
    "D65"  = -7.545
    "D69"  = -7.640000000000001
    "D79"  = -7.98
    "D83"  = -8.264999999999999
    "D91"  = -7.589
    "D93"  = -6.976999999999999
    "D100" = -7.969000000000001
    "D103" = -7.995999999999998
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
